$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.500.15"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "2.601.75"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.79"
$ws.Range("E5").Value = "  +2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.71"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").Value = "  +4.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.65"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.346"
$ws.Range("E11").Value = "  +2.06%  "

$ws.Range("E12").Value = "  +1.96%  "

$ws.Range("D13").Value = "3.055.38"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").Value = "60.528.93"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.61"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("E16").Value = "  +2.09%  "

$ws.Range("D17").Value = "2.601.19"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.75"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "358.48"
$ws.Range("E19").Value = "  +4.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.57"
$ws.Range("E20").Value = "  +3.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").Value = "  +3.27%  "

$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.03"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.427"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.716.08"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").Value = "0.0₃0837"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.44"
$ws.Range("E31").Value = "  +1.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.59"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.97"
$ws.Range("E33").Value = "  +5.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.48"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.02"
$ws.Range("E35").Value = "  +2.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.923"
$ws.Range("E36").Value = "  +9.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("E38").Value = "  +1.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.31"
$ws.Range("E39").Value = "  +2.57%  "

$ws.Range("E40").Value = "  +0.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "289.23"
$ws.Range("E42").Value = "  -1.74%  "

$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0556"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.58"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.96"
$ws.Range("E48").Value = "  +1.79%  "

$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.17"
$ws.Range("E51").Value = "  +11.11%  "
